$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This edit rotates several description paragraphs of the LOT2044 course
# sheet: each Portuguese "value" paragraph is replaced by the text that used
# to live in the paragraph that followed it (wrapping forward through the
# document), and the two now-redundant italic English paragraphs (the ones
# that sat right after "Objetivos" and right after "Programa resumido", both
# consisting of an empty/"Specific work plan..." run formatted with <w:i/>)
# are removed so that every remaining heading keeps exactly one body
# paragraph under it.
# ---------------------------------------------------------------------------

# 1) "Objetivos" body paragraph (was the PT objectives text) now holds the
#    text that used to be the "Programa resumido" PT paragraph.
$r = $d.Paragraphs(6).Range
$target = $d.Range($r.Start, $r.End - 1)
$target.Text = "Plano de Trabalho específico. Realização do Estágio. Relatório final e/ou parciais."

# 2) "Docente(s) Responsável(eis)" paragraph now holds the old "Objetivos"
#    PT text.
$r = $d.Paragraphs(9).Range
$target = $d.Range($r.Start, $r.End - 1)
$target.Text = "Fornecer oportunidade de aplicação dos conhecimentos fundamentais da Engenharia Bioquímica nos projetos e processos bioquímicos. Complementação da formação geral curricular. Adaptação psicológica e social do estudante à sua futura atividade profissional."

# 3) "Programa resumido" PT paragraph now holds the old "Programa" PT text.
$r = $d.Paragraphs(11).Range
$target = $d.Range($r.Start, $r.End - 1)
$target.Text = "Participação do aluno em processo seletivo de empresas ou no setor acadêmico. Estágio realizado sob a supervisão da Escola de Engenharia de Lorena, através do Departamento de Biotecnologia. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o Supervisor do Estágio e o professor orientador, desde que relacionado com as áreas afins da Engenharia Bioquímica.  Apresentação de relatório final e/ou relatórios parciais sobre as atividades desenvolvidas no estágio."

# 4) "Programa" PT paragraph now holds the old "Método" value text.
$r = $d.Paragraphs(14).Range
$target = $d.Range($r.Start, $r.End - 1)
$target.Text = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."

# 5) Inside the "Avaliação" bullet paragraph, the three labelled value runs
#    each take on the value that used to belong to the next label. Replace
#    right-to-left (last label first) so each Find.Execute only ever matches
#    the single still-original occurrence of its search text.
$r = $d.Paragraphs(17).Range
$r.Find.Execute("Não será oferecida recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "A ser definida com o orientador em função das atividades desenvolvidas no estágio.", 2) | Out-Null

$r = $d.Paragraphs(17).Range
$r.Find.Execute("MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio", $true, $false, $false, $false, $false, $true, 1, $false, "Não será oferecida recuperação.", 2) | Out-Null

$r = $d.Paragraphs(17).Range
$r.Find.Execute("Supervisão das atividades desenvolvidas pelo aluno durante o estágio.", $true, $false, $false, $false, $false, $true, 1, $false, "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio", 2) | Out-Null

# 6) "Bibliografia" paragraph now holds the old "Docente(s) Responsável(eis)"
#    text.
$r = $d.Paragraphs(19).Range
$target = $d.Range($r.Start, $r.End - 1)
$target.Text = "101761 - Arnaldo Márcio Ramalho Prata"

# 7) Remove the two now-orphaned italic English paragraphs, highest index
#    first so the lower index stays valid.
$d.Paragraphs(12).Range.Delete() | Out-Null   # italic "Specific work plan..." (Programa resumido EN)
$d.Paragraphs(7).Range.Delete() | Out-Null    # empty italic run (Objetivos EN)
